$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" list
# ("LOM3234: ...") so we can remove the three paragraphs that
# used to follow it (a blank paragraph, the "Ver no Jupiter ..."
# line and the "(c) 2020 ..." footer line) while leaving the
# trailing blank paragraph (and the page-break paragraph after it)
# untouched.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*LOM3234*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $start = $d.Paragraphs.Item($anchorIndex + 1).Range.Start
    $end = $d.Paragraphs.Item($anchorIndex + 3).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
